# Update the answers in the "two-digit ÷ one-digit" practice table.
# Each cell is addressed directly by (row, column) in the single table
# so that the replacement is unambiguous even though several of the
# new values coincide textually with old values used elsewhere in the
# sheet (which would make a naive global Find/Replace unsafe).

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "32÷9=3, 5" },
    @{ Row = 1;  Col = 2; Text = "81÷2=40, 1" },
    @{ Row = 1;  Col = 3; Text = "23÷7=3, 2" },
    @{ Row = 1;  Col = 4; Text = "65÷8=8, 1" },
    @{ Row = 1;  Col = 5; Text = "23÷4=5, 3" },

    @{ Row = 5;  Col = 1; Text = "54÷3=18, 0" },
    @{ Row = 5;  Col = 2; Text = "88÷9=9, 7" },
    @{ Row = 5;  Col = 3; Text = "45÷2=22, 1" },
    @{ Row = 5;  Col = 4; Text = "34÷9=3, 7" },
    @{ Row = 5;  Col = 5; Text = "66÷2=33, 0" },

    @{ Row = 9;  Col = 1; Text = "82÷9=9, 1" },
    @{ Row = 9;  Col = 2; Text = "21÷6=3, 3" },
    @{ Row = 9;  Col = 3; Text = "33÷9=3, 6" },
    @{ Row = 9;  Col = 4; Text = "90÷7=12, 6" },
    @{ Row = 9;  Col = 5; Text = "75÷5=15, 0" },

    @{ Row = 13; Col = 1; Text = "96÷9=10, 6" },
    @{ Row = 13; Col = 2; Text = "34÷4=8, 2" },
    @{ Row = 13; Col = 3; Text = "32÷4=8, 0" },
    @{ Row = 13; Col = 4; Text = "36÷5=7, 1" },
    @{ Row = 13; Col = 5; Text = "82÷8=10, 2" },

    @{ Row = 17; Col = 1; Text = "68÷8=8, 4" },
    @{ Row = 17; Col = 2; Text = "73÷3=24, 1" },
    @{ Row = 17; Col = 3; Text = "88÷8=11, 0" },
    @{ Row = 17; Col = 4; Text = "34÷2=17, 0" },
    @{ Row = 17; Col = 5; Text = "32÷8=4, 0" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
